$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Params")

# --- New row 16: coefficient of O2 solubility in water (Henry's law) ---
# Write in the order that reproduces the target shared-string table order.
$ws.Range("A16").Value = "coefficient of O2 solubility in water"
$ws.Range("B16").Formula = "=1.39*10^-3"
$ws.Range("C16").Value = "mMol*L^-1 *mmHg^-1"

# --- Row 13 (K_M_RNR): convert from mmHg to mMol*L^-1 using the new
#     solubility coefficient in B16 ---
$ws.Range("C13").Value = "mMol*L^-1"
$ws.Range("B13").Formula = "=0.009241*B15*B16"

# Copy the number/text formatting (style) from the analogous cells in row 15
# onto the new row 16 cells, preserving the values already written.
$ws.Range("D15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = "Henry's law"

$ws.Range("D15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "12771568"

$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)

# --- Cosmetic sheet-view changes ---
$ws.Columns.Item(1).ColumnWidth = 37.5
$ws.Range("D16").Select()
